# Apply updated dSF (column F) values for specific rows, per repull of data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 3
    3  = 1
    4  = 0
    8  = -1
    13 = -3
    16 = 7
    19 = -2
    20 = 1
    23 = -1
    25 = 0
    26 = 6
    27 = -2
    31 = 6
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
